$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 93.5625
$ws.Range("I15").Value = 93.5625
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 280.6875
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = -111.6875

$ws.Range("H64").Value = 3799.8
$ws.Range("I64").Value = 3499.5
$ws.Range("J64").Value = 4000
$ws.Range("K64").Value = 3499.5
$ws.Range("L64").Value = 4000
$ws.Range("M64").Value = -3251.5
$ws.Range("N64").Value = -4496

$ws.Range("H67").Value = 3799.8
$ws.Range("I67").Value = 3499.5
$ws.Range("J67").Value = 4000
$ws.Range("K67").Value = 3499.5
$ws.Range("L67").Value = 4000
$ws.Range("M67").Value = -2641.5
$ws.Range("N67").Value = -5716

$ws.Range("H100").Value = 4715.8335
$ws.Range("I100").Value = 4598.3335
$ws.Range("J100").Value = 4833.3335
$ws.Range("K100").Value = 4598.3335
$ws.Range("L100").Value = 4833.3335
$ws.Range("M100").Value = -4057.3335
$ws.Range("N100").Value = -5915.3335

$ws.Range("H103").Value = 613.6923
$ws.Range("I103").Value = 288
$ws.Range("J103").Value = 892.8570999999999
$ws.Range("K103").Value = 864
$ws.Range("L103").Value = 2678.5713
$ws.Range("M103").Value = -278
$ws.Range("N103").Value = -3850.5713

$ws.Range("H125").Value = 2507.75
$ws.Range("I125").Value = 2331.6667
$ws.Range("J125").Value = 3036
$ws.Range("K125").Value = 20985.0003
$ws.Range("L125").Value = 27324
$ws.Range("M125").Value = -18525.0003
$ws.Range("N125").Value = -32244

$ws.Range("H127").Value = 741
$ws.Range("I127").Value = 741
$ws.Range("J127").Value = 0
$ws.Range("K127").Value = 2223
$ws.Range("L127").Value = 0
$ws.Range("M127").Value = 2737

$ws.Range("H131").Value = 2557.6
$ws.Range("I131").Value = 798.3333
$ws.Range("J131").Value = 5196.5
$ws.Range("K131").Value = 2394.9999
$ws.Range("L131").Value = 15589.5
$ws.Range("M131").Value = 2645.0001
$ws.Range("N131").Value = -25669.5

$ws.Range("H132").Value = 4098.1665
$ws.Range("I132").Value = 2124.2727
$ws.Range("J132").Value = 7200
$ws.Range("K132").Value = 6372.8181
$ws.Range("L132").Value = 21600
$ws.Range("M132").Value = -3842.8181
$ws.Range("N132").Value = -26660

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 6375
$ws.Range("I2").Value = 5750
$ws.Range("J2").Value = 7000
$ws.Range("K2").Value = 5750
$ws.Range("L2").Value = 7000
$ws.Range("M2").Value = -5637
$ws.Range("N2").Value = -7226

$ws.Range("H32").Value = 4185.359
$ws.Range("I32").Value = 3098.5
$ws.Range("J32").Value = 11576
$ws.Range("K32").Value = 3098.5
$ws.Range("L32").Value = 11576
$ws.Range("M32").Value = -2811.5
$ws.Range("N32").Value = -12150

$ws.Range("H116").Value = 6375
$ws.Range("I116").Value = 5750
$ws.Range("J116").Value = 7000
$ws.Range("K116").Value = 5750
$ws.Range("L116").Value = 7000
$ws.Range("M116").Value = -3456
$ws.Range("N116").Value = -11588

$ws.Range("H132").Value = 1706.5294
$ws.Range("I132").Value = 1706.5294
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 5119.5882
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -2589.5882

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 6375
$ws.Range("I3").Value = 5750
$ws.Range("J3").Value = 7000
$ws.Range("K3").Value = 5750
$ws.Range("L3").Value = 7000
$ws.Range("M3").Value = -5636
$ws.Range("N3").Value = -7228

$ws.Range("H80").Value = 674.7778
$ws.Range("I80").Value = 788.5
$ws.Range("J80").Value = 583.8
$ws.Range("K80").Value = 788.5
$ws.Range("L80").Value = 583.8
$ws.Range("M80").Value = 209.5
$ws.Range("N80").Value = -2579.8

$ws.Range("H83").Value = 674.7778
$ws.Range("I83").Value = 788.5
$ws.Range("J83").Value = 583.8
$ws.Range("K83").Value = 3942.5
$ws.Range("L83").Value = 2919
$ws.Range("M83").Value = 1049.5
$ws.Range("N83").Value = -12903

$ws.Range("H134").Value = 3671.3333
$ws.Range("I134").Value = 4500
$ws.Range("J134").Value = 2014
$ws.Range("K134").Value = 13500
$ws.Range("L134").Value = 6042
$ws.Range("M134").Value = -10965
$ws.Range("N134").Value = -11112

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 7000600
$ws.Range("I6").Value = 3333667
$ws.Range("J6").Value = 12501000
$ws.Range("K6").Value = 3333667
$ws.Range("L6").Value = 12501000
$ws.Range("M6").Value = -3333554
$ws.Range("N6").Value = -12501226

$ws.Range("H31").Value = 9498.6
$ws.Range("I31").Value = 12497.667
$ws.Range("J31").Value = 5000
$ws.Range("K31").Value = 12497.667
$ws.Range("L31").Value = 5000
$ws.Range("M31").Value = -12202.667
$ws.Range("N31").Value = -5590

$ws.Range("H34").Value = 9498.6
$ws.Range("I34").Value = 12497.667
$ws.Range("J34").Value = 5000
$ws.Range("K34").Value = 12497.667
$ws.Range("L34").Value = 5000
$ws.Range("M34").Value = -12295.667
$ws.Range("N34").Value = -5404

$ws.Range("H62").Value = 7709
$ws.Range("I62").Value = 8444.875
$ws.Range("J62").Value = 6727.8335
$ws.Range("K62").Value = 8444.875
$ws.Range("L62").Value = 6727.8335
$ws.Range("M62").Value = -7820.875
$ws.Range("N62").Value = -7975.8335

$ws.Range("H65").Value = 7709
$ws.Range("I65").Value = 8444.875
$ws.Range("J65").Value = 6727.8335
$ws.Range("K65").Value = 42224.375
$ws.Range("L65").Value = 33639.1675
$ws.Range("M65").Value = -39104.375
$ws.Range("N65").Value = -39879.1675

$ws.Range("H99").Value = 1497
$ws.Range("I99").Value = 1497
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 1497
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = 1

$ws.Range("H126").Value = 1497
$ws.Range("I126").Value = 1497
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 4491
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -2021

$ws.Range("H132").Value = 3141
$ws.Range("I132").Value = 2863.75
$ws.Range("J132").Value = 4250
$ws.Range("K132").Value = 8591.25
$ws.Range("L132").Value = 12750
$ws.Range("M132").Value = -6061.25
$ws.Range("N132").Value = -17810

$ws.Range("H134").Value = 2750
$ws.Range("I134").Value = 2750
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 8250
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -5715

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("M16").ClearContents()

$ws.Range("H68").Value = 1696.75
$ws.Range("I68").Value = 2000
$ws.Range("J68").Value = 1595.6666
$ws.Range("K68").Value = 6000
$ws.Range("L68").Value = 4786.9998
$ws.Range("M68").Value = -5189
$ws.Range("N68").Value = -6408.9998

$ws.Range("H71").Value = 1696.75
$ws.Range("I71").Value = 2000
$ws.Range("J71").Value = 1595.6666
$ws.Range("K71").Value = 18000
$ws.Range("L71").Value = 14360.9994
$ws.Range("M71").Value = -13944
$ws.Range("N71").Value = -22472.9994

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3137.5334
$ws.Range("I132").Value = 3428
$ws.Range("J132").Value = 1249.5
$ws.Range("K132").Value = 10284
$ws.Range("L132").Value = 3748.5
$ws.Range("M132").Value = -7754
$ws.Range("N132").Value = -8808.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7055.5557
$ws.Range("I7").Value = 6833.3335
$ws.Range("J7").Value = 7500
$ws.Range("K7").Value = 6833.3335
$ws.Range("L7").Value = 7500
$ws.Range("M7").Value = -6721.3335
$ws.Range("N7").Value = -7724

$ws.Range("H126").Value = 7055.5557
$ws.Range("I126").Value = 6833.3335
$ws.Range("J126").Value = 7500
$ws.Range("K126").Value = 20500.0005
$ws.Range("L126").Value = 22500
$ws.Range("M126").Value = -18030.0005
$ws.Range("N126").Value = -27440

$ws.Range("H132").Value = 6562.125
$ws.Range("I132").Value = 6562.125
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 19686.375
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -17156.375

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 14191
$ws.Range("I100").Value = 14191
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 28382
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -27841

$ws.Range("H122").Value = 2898.5454
$ws.Range("I122").Value = 2554.4285
$ws.Range("J122").Value = 3500.75
$ws.Range("K122").Value = 7663.2855
$ws.Range("L122").Value = 10502.25
$ws.Range("M122").Value = -5213.2855
$ws.Range("N122").Value = -15402.25

$ws.Range("H132").Value = 2348.6667
$ws.Range("I132").Value = 2418.4
$ws.Range("J132").Value = 2000
$ws.Range("K132").Value = 7255.200000000001
$ws.Range("L132").Value = 6000
$ws.Range("M132").Value = -4725.200000000001
$ws.Range("N132").Value = -11060

$ws.Range("H136").Value = 2907.5715
$ws.Range("I136").Value = 2907.5715
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 8722.7145
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -6172.7145
$ws.Range("N136").ClearContents()
